$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H (vc-espectaculos): was curated as a dimension, now curated as a measure.
$ws.Range("H2").Value = "iaest-measure:vc-espectaculos"
$ws.Range("H3").Value = "medida"
$ws.Range("H4").Value = "xsd:int"
$ws.Range("H5").Clear()

# Column U (municipio-nombre): was curated as a measure, now curated as a dimension.
$ws.Range("U2").Value = "sdmx-dimension:refArea"
$ws.Range("U3").Value = "dim"
$ws.Range("U4").Value = "URI-Municipio"
